$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text correction: "Limit Switches" -> "Limit and Home Switches" ---
$ws.Range("A10").Value2 = "Limit and Home Switches"

# --- New Description (column F) entries for the CNC Machine section ---
$ws.Range("F6").Value2  = "Stepper motor for the X axis"
$ws.Range("F7").Value2  = "Stepper Motor for the Y axis"
$ws.Range("F8").Value2  = "Stepper Moto for the Z axis"
$ws.Range("F9").Value2  = "PCI I/O Card with Digital and Analog I/O"
$ws.Range("F11").Value2 = "This is part of the CNC System"
$ws.Range("F12").Value2 = "Interface Card for communicating with servo motors"
$ws.Range("F13").Value2 = "I/O card used for limit and home switches"
$ws.Range("F14").Value2 = "Temp range x -x "
$ws.Range("F15").Value2 = " Outputs a single Voltage"
$ws.Range("F16").Value2 = "Wire Feed Welder that uses inert sheilding gas"
$ws.Range("F17").Value2 = "Used for controlling the motors on the welder control knobs"
$ws.Range("F18").Value2 = "Used for connecting signals to the PCIe card"

# --- New quantity cells for the two STL-printed assemblies ---
$ws.Range("E19").Value2 = 5
$ws.Range("E23").Value2 = 5

# --- Widen column F to fit the new Description text ---
$ws.Columns.Item(6).ColumnWidth = 49.4986979166667

# --- Move the active selection (as recorded in the saved view state) ---
$ws.Range("E21").Select() | Out-Null
